# Edit script: add a "type" classification column (K) to the hidden
# "grammar_1" sheet, and update the saved selections on both sheets, to
# match the upstream commit that added column K with grammar-topic labels.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("grammar_1")
$ws2 = $wb.Worksheets.Item("Sheet1")

# --- New column K: "type" header + grammar-topic labels per question block ---
$ws1.Range("K1").Value = "type"

$ws1.Range("K2:K5").Value   = "명사의 복수형"
$ws1.Range("K6:K9").Value   = "물질명사와 단위명사"
$ws1.Range("K10:K16").Value = "관사"
$ws1.Range("K17:K20").Value = "지시대명사와 지시형용사"
$ws1.Range("K21:K28").Value = "인칭대명사"
$ws1.Range("K29:K36").Value = "be동사"
$ws1.Range("K37:K41").Value = "be동사의 평사문, 부정문, 의문문"
$ws1.Range("K42:K46").Value = "일반동사"
$ws1.Range("K47:K51").Value = "의문문"

# --- Update the saved cursor/selection & scroll position on grammar_1 ---
$ws1.Activate()
$excel.ActiveWindow.ScrollRow    = 31
$excel.ActiveWindow.ScrollColumn = 2
$ws1.Range("I39").Select()

# --- Update the saved cursor/selection on the visible Sheet1 (stays active) ---
$ws2.Activate()
$ws2.Range("C25").Select()
